$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("aryl_halide_200epochs","ranking_test2",5,0.001,"mean","sum",212.385,0.725,14.202,0.521,16.653),
    @("aryl_halide_200epochs","ranking_test2",4,0.01,"max","sum",216.155,0.725,14.233,0.5580000000000001,16.23),
    @("aryl_halide_200epochs","ranking_test2",4,0.001,"max","sum",218.584,0.724,14.226,0.542,16.425),
    @("aryl_halide_200epochs","ranking_test2",5,0.01,"max","sum",226.364,0.718,14.418,0.54,16.466),
    @("aryl_halide_200epochs","ranking_test2",5,0.001,"max","sum",221.295,0.725,14.229,0.556,16.217),
    @("aryl_halide_200epochs","ranking_test2",4,0.001,"mean","sum",223.117,0.715,14.447,0.487,17.015),
    @("aryl_halide_200epochs","ranking_test2",3,0.01,"max","sum",233.778,0.707,14.642,0.533,16.655),
    @("aryl_halide_200epochs","ranking_test2",5,0.01,"max","mean",221.283,0.724,14.282,0.543,16.392),
    @("aryl_halide_200epochs","ranking_test2",4,0.01,"max","product",211.625,0.76,13.205,0.552,15.916),
    @("aryl_halide_200epochs","ranking_test2",5,0.001,"mean","product",160.118,0.794,12.237,0.611,15.004),
    @("aryl_halide_200epochs","ranking_test2",4,0.01,"mean","mean",221.944,0.725,14.261,0.392,18.457),
    @("aryl_halide_200epochs","ranking_test2",5,0.0001,"mean","sum",258.571,0.662,15.756,0.406,18.23),
    @("aryl_halide_200epochs","ranking_test2",4,0.0001,"max","sum",281.136,0.627,16.553,0.271,20.184),
    @("aryl_halide_200epochs","ranking_test2",5,0.0001,"max","sum",270.144,0.645,16.129,0.272,20.26),
    @("aryl_halide_200epochs","ranking_test2",4,0.0001,"mean","sum",293.201,0.615,16.763,0.389,18.645),
    @("aryl_halide_200epochs","ranking_test2",3,0.001,"max","sum",226.792,0.704,14.697,0.518,16.866),
    @("aryl_halide_200epochs","ranking_test2",3,0.001,"max","sum",235.59,0.6899999999999999,15.104,0.504,16.966),
    @("aryl_halide_200epochs","ranking_test2",5,0.001,"max","mean",216.875,0.724,14.203,0.573,16.022),
    @("aryl_halide_200epochs","ranking_test2",5,0.0001,"max","mean",267.155,0.647,16.138,0.333,19.286),
    @("aryl_halide_200epochs","ranking_test2",4,0.001,"max","product",160.122,0.795,12.255,0.654,14.14),
    @("aryl_halide_200epochs","ranking_test2",4,0.0001,"max","product",195.057,0.754,13.367,0.5659999999999999,15.917),
    @("aryl_halide_200epochs","ranking_test2",5,0.0001,"mean","product",207.858,0.744,13.726,0.517,16.753),
    @("aryl_halide_200epochs","ranking_test2",4,0.001,"mean","mean",230.094,0.703,14.678,0.501,16.957),
    @("aryl_halide_200epochs","ranking_test2",4,0.0001,"mean","mean",320.148,0.581,17.47,0.352,19.092)
)

$startRow = 74
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    for ($col = 1; $col -le 11; $col++) {
        $ws.Cells.Item($row, $col).Value = $rec[$col - 1]
    }
}
